{"js": "// The document's paragraphs (top-level body paragraph + every table-cell\n// paragraph, in document order) map 1:1 onto the text runs touched by the\n// diff: index 0 is the date line, and each \"problem row\" of the practice\n// table occupies 5 consecutive paragraph slots separated by 3 blank-answer\n// rows. Replacing by position (rather than a global text search) is\n// required because several source strings (e.g. \"22\u00f73=\") repeat at\n// different positions with different replacements.\nconst replacements = {\n  0: { from: \"2025-12-18 Thursday\", to: \"2025-12-19 Friday\" },\n  1: { from: \"13\u00f79=\", to: \"88\u00f78=\" },\n  2: { from: \"55\u00f78=\", to: \"76\u00f79=\" },\n  3: { from: \"70\u00f77=\", to: \"73\u00f77=\" },\n  4: { from: \"90\u00f73=\", to: \"50\u00f76=\" },\n  5: { from: \"22\u00f73=\", to: \"67\u00f77=\" },\n  21: { from: \"22\u00f73=\", to: \"96\u00f77=\" },\n  22: { from: \"47\u00f79=\", to: \"46\u00f77=\" },\n  23: { from: \"86\u00f75=\", to: \"63\u00f74=\" },\n  24: { from: \"23\u00f76=\", to: \"96\u00f76=\" },\n  25: { from: \"74\u00f72=\", to: \"19\u00f73=\" },\n  41: { from: \"33\u00f74=\", to: \"92\u00f74=\" },\n  42: { from: \"38\u00f77=\", to: \"81\u00f79=\" },\n  43: { from: \"52\u00f72=\", to: \"90\u00f75=\" },\n  44: { from: \"30\u00f77=\", to: \"91\u00f73=\" },\n  45: { from: \"94\u00f79=\", to: \"11\u00f72=\" },\n  61: { from: \"10\u00f75=\", to: \"14\u00f73=\" },\n  62: { from: \"46\u00f73=\", to: \"59\u00f76=\" },\n  63: { from: \"86\u00f77=\", to: \"38\u00f77=\" },\n  64: { from: \"92\u00f75=\", to: \"56\u00f78=\" },\n  65: { from: \"66\u00f72=\", to: \"93\u00f75=\" },\n  81: { from: \"73\u00f78=\", to: \"88\u00f73=\" },\n  82: { from: \"29\u00f76=\", to: \"57\u00f72=\" },\n  83: { from: \"46\u00f75=\", to: \"15\u00f79=\" },\n  84: { from: \"72\u00f72=\", to: \"37\u00f73=\" },\n  85: { from: \"64\u00f79=\", to: \"57\u00f74=\" },\n};\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst maxIndex = Math.max(...Object.keys(replacements).map(Number));\nif (paragraphs.items.length <= maxIndex) {\n  throw new Error(\n    `Expected at least ${maxIndex + 1} paragraphs, found ${paragraphs.items.length}`\n  );\n}\n\n// Verify the paragraphs still hold the text the diff was generated against\n// before mutating anything, so a structurally different document fails\n// loudly instead of silently overwriting the wrong cell.\nfor (const [indexStr, { from }] of Object.entries(replacements)) {\n  const index = Number(indexStr);\n  const actual = paragraphs.items[index].text;\n  if (actual !== from) {\n    throw new Error(\n      `Paragraph ${index}: expected ${JSON.stringify(from)}, found ${JSON.stringify(actual)}`\n    );\n  }\n}\n\nfor (const [indexStr, { to }] of Object.entries(replacements)) {\n  const index = Number(indexStr);\n  paragraphs.items[index].insertText(to, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# The practice sheet's date line becomes the next day, and every division\n# problem in the table is swapped for a new one. Several source problems\n# (e.g. \"22\u00f73=\") repeat at different grid positions with different\n# replacements, so cells are addressed positionally via Table.Cell(row,col)\n# rather than via a global text search/replace.\n\n$d = $word.ActiveDocument\n\n# Date heading (first paragraph in the document body).\n$d.Paragraphs.Item(1).Range.Text = \"2025-12-19 Friday\"\n\n$t = $d.Tables.Item(1)\n\n# Map of 1-indexed (row, column) -> new problem text. Rows 1,5,9,13,17 hold\n# the five problems; the three rows between each are blank answer rows.\n$updates = @{\n  \"1,1\" = \"88\u00f78=\"\n  \"1,2\" = \"76\u00f79=\"\n  \"1,3\" = \"73\u00f77=\"\n  \"1,4\" = \"50\u00f76=\"\n  \"1,5\" = \"67\u00f77=\"\n\n  \"5,1\" = \"96\u00f77=\"\n  \"5,2\" = \"46\u00f77=\"\n  \"5,3\" = \"63\u00f74=\"\n  \"5,4\" = \"96\u00f76=\"\n  \"5,5\" = \"19\u00f73=\"\n\n  \"9,1\" = \"92\u00f74=\"\n  \"9,2\" = \"81\u00f79=\"\n  \"9,3\" = \"90\u00f75=\"\n  \"9,4\" = \"91\u00f73=\"\n  \"9,5\" = \"11\u00f72=\"\n\n  \"13,1\" = \"14\u00f73=\"\n  \"13,2\" = \"59\u00f76=\"\n  \"13,3\" = \"38\u00f77=\"\n  \"13,4\" = \"56\u00f78=\"\n  \"13,5\" = \"93\u00f75=\"\n\n  \"17,1\" = \"88\u00f73=\"\n  \"17,2\" = \"57\u00f72=\"\n  \"17,3\" = \"15\u00f79=\"\n  \"17,4\" = \"37\u00f73=\"\n  \"17,5\" = \"57\u00f74=\"\n}\n\nforeach ($row in @(1, 5, 9, 13, 17)) {\n  for ($col = 1; $col -le 5; $col++) {\n    $key = \"$row,$col\"\n    $t.Cell($row, $col).Range.Text = $updates[$key]\n  }\n}\n"}
